# Apply the "stuff at the bottom of the sheets" commit:
#  - fill in the (previously blank) pair_kind column for the practice rows
#  - append a new "stim details" block describing month/word_type/need_audio/
#    need_image/word/count/find images, with sample rows for video & audio

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows (2-5) were missing a pair_kind (column J) value; they are
# all "generic" pairs, same as the rest of the sheet.
$ws.Range("J2:J5").Value = "generic"

# New block starting at row 27.
$ws.Range("A27").Value = "stim details"

# Header row for the new block.
$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# Data rows for the new block.
$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "video"

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "video"

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "video"

$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "video"

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "audio"

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "audio"

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "audio"

$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "audio"
